$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "V+" / "V-" header labels in B1 and C1
$ws.Range("B1").Value = "V-"
$ws.Range("C1").Value = "V+"

# Move the active selection to C2 (reflects the cell last interacted with)
$ws.Range("C2").Select()
